$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "...в одном направление, а его выбор..." -> "...в одном направлении, а его выбор..."
$rng = $d.Content
$ok1 = $rng.Find.Execute(
    "в одном направление, а его выбор",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "в одном направлении, а его выбор",
    2)

# --- Change 2 -------------------------------------------------------------
# "), напряжение затвора имеет допустимое значение 3.3В  (подключим" ->
# "), напряжение затвора 3.3В для которого является допустимым (подключим"
$rng2 = $d.Content
$ok2 = $rng2.Find.Execute(
    "), напряжение затвора имеет допустимое значение 3.3В  (подключим",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "), напряжение затвора 3.3В для которого является допустимым (подключим",
    2)

# --- Change 3 -------------------------------------------------------------
# "ток сток-исток 0.825мА" -> "ток сток-исток составляет 0.825мА"
$rng3 = $d.Content
$ok3 = $rng3.Find.Execute(
    "ток сток-исток 0.825мА",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ток сток-исток составляет 0.825мА",
    2)

# --- Change 4 -------------------------------------------------------------
# "превышает в 1.5 раза" -> "превышает в 1.65 раза"
$rng4 = $d.Content
$ok4 = $rng4.Find.Execute(
    "превышает в 1.5 раза",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "превышает в 1.65 раза",
    2)

Write-Output "change1=$ok1 change2=$ok2 change3=$ok3 change4=$ok4"

if (-not ($ok1 -and $ok2 -and $ok3 -and $ok4)) {
    throw "One or more Find/Replace operations failed to locate their target text."
}
